# Applies cryptos list price/volume updates generated from the GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric ("582.79", "0.999", ...) must stay plain text
# (matches the source data which stores every Price/Volume cell as a string).
# Force text format first so Excel does not coerce the assigned string into a number,
# then clear the formatting change afterwards so the cell keeps the workbook default style.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "60.452.51"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "2.614.14"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue $ws.Range("D5") "582.79"
$ws.Range("E5").Value = "  +1.88%  "
Set-TextValue $ws.Range("D6") "143.59"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").Value = "3.074.36"
$ws.Range("E13").Value = "  -0.01%  "
Set-TextValue $ws.Range("D14") "24.71"
$ws.Range("E14").Value = "  +5.18%  "
$ws.Range("D15").Value = "60.414.45"
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").Value = "2.617.39"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("E19").Value = "  -0.97%  "
Set-TextValue $ws.Range("D20") "346.88"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("E21").Value = "  -2.95%  "
Set-TextValue $ws.Range("D22") "0.999"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("E23").Value = "  +1.36%  "
Set-TextValue $ws.Range("D24") "63.63"
$ws.Range("E24").Value = "  -0.66%  "
Set-TextValue $ws.Range("D25") "0.999"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  -0.42%  "
Set-TextValue $ws.Range("D27") "7.98"
$ws.Range("E27").Value = "  +3.24%  "
$ws.Range("E28").Value = "  +5.27%  "
$ws.Range("D29").Value = "0.0₃0799"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  +1.79%  "
Set-TextValue $ws.Range("D31") "168.76"
$ws.Range("E31").Value = "  +4.82%  "
$ws.Range("E32").Value = "  +0.18%  "
Set-TextValue $ws.Range("D33") "19.47"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("E34").Value = "  +8.60%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D35") "1.01"
$ws.Range("E35").Value = "  +4.49%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D36") "4.29"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  +3.55%  "
Set-TextValue $ws.Range("D38") "319.97"
$ws.Range("E38").Value = "  +7.17%  "
Set-TextValue $ws.Range("D39") "38.36"
$ws.Range("E40").Value = "  +2.38%  "
$ws.Range("E41").Value = "  -0.67%  "
Set-TextValue $ws.Range("D42") "136.02"
$ws.Range("E42").Value = "  -2.75%  "
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("E45").Value = "  +0.89%  "
Set-TextValue $ws.Range("D46") "19.95"
$ws.Range("E46").Value = "  +1.53%  "
Set-TextValue $ws.Range("D48") "0.0551"
$ws.Range("E48").Value = "  -0.09%  "
Set-TextValue $ws.Range("D49") "20.09"
$ws.Range("E49").Value = "  +1.36%  "
$ws.Range("E50").Value = "  -0.35%  "
Set-TextValue $ws.Range("D51") "10.76"
$ws.Range("E51").Value = "  +0.55%  "
